$d = $word.ActiveDocument

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading paragraph so we only touch
# the bullets that belong to that section (the document also contains a very
# similarly-worded bullet list under PROFESSIONAL EXPERIENCE that must stay
# untouched).
$headingIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*KEY ACHIEVEMENTS AND IMPACT*") {
        $headingIndex = $idx
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find KEY ACHIEVEMENTS AND IMPACT heading"
}

# Layout right after the heading:
#   headingIndex + 0 -> "KEY ACHIEVEMENTS AND IMPACT"
#   headingIndex + 1 -> "Impact"
#   headingIndex + 2 -> bullet 1 (Discovered systematic race coding errors...)
#   headingIndex + 3 -> bullet 2 (Trigonometric algorithm...)
#   headingIndex + 4 -> bullet 3 (Built redistricting platform...)
#   headingIndex + 5 -> bullet 4 (Achieved 87% prediction accuracy...)
#   headingIndex + 6 -> bullet 5 (Built real-time FEC analysis systems...)  [removed]
#   headingIndex + 7 -> bullet 6 (Provided expert testimony...)

$bullet1 = $d.Paragraphs.Item($headingIndex + 2)
$bullet2 = $d.Paragraphs.Item($headingIndex + 3)
$bullet3 = $d.Paragraphs.Item($headingIndex + 4)
$bullet4 = $d.Paragraphs.Item($headingIndex + 5)
$bullet6 = $d.Paragraphs.Item($headingIndex + 7)

$bullet1.Range.Find.Execute("Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%", $true, $false, $false, $false, $false, $true, 1, $false, "Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions", 2) | Out-Null

$bullet2.Range.Find.Execute("Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis", $true, $false, $false, $false, $false, $true, 1, $false, "178% accuracy improvement in racial classification algorithms", 2) | Out-Null

$bullet3.Range.Find.Execute("Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations", $true, $false, $false, $false, $false, $true, 1, $false, "Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%", 2) | Out-Null

$bullet4.Range.Find.Execute("Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%", $true, $false, $false, $false, $false, $true, 1, $false, "$4.7M savings enabled nonprofit access", 2) | Out-Null

$bullet6.Range.Find.Execute("Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy", $true, $false, $false, $false, $false, $true, 1, $false, "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations", 2) | Out-Null

# Remove the "Built real-time FEC analysis systems..." bullet entirely.
$bullet5 = $d.Paragraphs.Item($headingIndex + 6)
$bullet5.Range.Delete()
